$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.05840966666666667
$ws.Range("H2").Value = 0.175229

$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.05840966666666667
$ws.Range("N2").Value = 0.175229

$ws.Range("Q2").Value = 0.003411689160111112
$ws.Range("R2").Value = 0.03070520244100001
